# "remove column from alcohol data"
# The alcohol measurement sheet (Sheet1) has an extra column M that duplicates/
# shifts into N. Delete column M entirely so N's data shifts left into M.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("M:M").EntireColumn.Delete()

# Matches the post-edit selection recorded in the workbook (M1 on Sheet1).
$ws.Range("M1").Select()
